$d = $word.ActiveDocument

function Find-ParaIndex($doc, $text) {
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        if ($doc.Paragraphs.Item($i).Range.Text -eq ($text + "`r")) {
            return $i
        }
    }
    return -1
}

$pkgOpen = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>'
$pkgClose = '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

# -----------------------------------------------------------------
# Edit 1: insert a new paragraph "The Other Side" right before the
# "Resources" paragraph.
# -----------------------------------------------------------------
$idx = Find-ParaIndex $d "Resources"
$insertRange = $d.Paragraphs.Item($idx).Range
$insertRange.Collapse(1)  # wdCollapseStart
$null = $insertRange.InsertParagraphBefore()

$theOtherSideXml = $pkgOpen + '<w:p><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:ind w:left="360"/><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="000000"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial" w:eastAsia="Times New Roman"/><w:color w:val="000000"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:tab/><w:t xml:space="preserve">The Other Side</w:t></w:r></w:p>' + $pkgClose
$d.Paragraphs.Item($idx).Range.InsertXML($theOtherSideXml)

# -----------------------------------------------------------------
# Edit 2: before the "Graphics" paragraph, insert two new paragraphs:
#   "Exchanges Listing Guide" and a new "Graphics" paragraph; then
#   change the original "Graphics" paragraph's run to "Marketing
#   Materials".
# -----------------------------------------------------------------
$idx = Find-ParaIndex $d "Graphics"

$r1 = $d.Paragraphs.Item($idx).Range
$r1.Collapse(1)
$null = $r1.InsertParagraphBefore()

$r2 = $d.Paragraphs.Item($idx + 1).Range
$r2.Collapse(1)
$null = $r2.InsertParagraphBefore()

$exchangesGuideXml = $pkgOpen + '<w:p><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:ind w:left="360"/><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="000000"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial" w:eastAsia="Times New Roman"/><w:color w:val="000000"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:tab/><w:tab/><w:t xml:space="preserve">Exchanges Listing Guide</w:t></w:r></w:p>' + $pkgClose
$d.Paragraphs.Item($idx).Range.InsertXML($exchangesGuideXml)

$newGraphicsXml = $pkgOpen + '<w:p><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:ind w:left="360" w:firstLine="360"/><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="000000"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial" w:eastAsia="Times New Roman"/><w:color w:val="000000"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>Graphics</w:t></w:r></w:p>' + $pkgClose
$d.Paragraphs.Item($idx + 1).Range.InsertXML($newGraphicsXml)

# The original "Graphics" paragraph is now at $idx + 2; replace its
# run (but keep the paragraph itself / its pPr) with "Marketing
# Materials".
$origGraphicsPara = $d.Paragraphs.Item($idx + 2)
$fullRange = $origGraphicsPara.Range
$runOnlyRange = $d.Range($fullRange.Start, $fullRange.End - 1)
$marketingMaterialsXml = $pkgOpen + '<w:p><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman" w:eastAsia="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:tab/><w:t xml:space="preserve">Marketing Materials</w:t></w:r></w:p>' + $pkgClose
$runOnlyRange.InsertXML($marketingMaterialsXml)

# -----------------------------------------------------------------
# Edit 3: insert a new paragraph "Guide" right after the "Electrum
# Wallet" paragraph.
# -----------------------------------------------------------------
$idx = Find-ParaIndex $d "Electrum Wallet"
$afterRange = $d.Paragraphs.Item($idx).Range
$afterRange.Collapse(0)  # wdCollapseEnd
$null = $afterRange.InsertParagraphAfter()

$guideXml = $pkgOpen + '<w:p><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:ind w:left="360" w:firstLine="360"/><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="000000"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial" w:eastAsia="Times New Roman"/><w:color w:val="000000"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/></w:rPr><w:tab/><w:t>Guide</w:t></w:r></w:p>' + $pkgClose
$d.Paragraphs.Item($idx + 1).Range.InsertXML($guideXml)

Write-Host "All edits applied."
